$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=5; I='sv'; J='Statement-opinion'},
    @{Row=6; I='ba'; J='Appreciation'},
    @{Row=25; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=38; I='sv'; J='Statement-opinion'},
    @{Row=50; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=59; I='sv'; J='Statement-opinion'},
    @{Row=67; I='ba'; J='Appreciation'},
    @{Row=87; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=95; I='sv'; J='Statement-opinion'},
    @{Row=96; I='sv'; J='Statement-opinion'},
    @{Row=102; I='sv'; J='Statement-opinion'},
    @{Row=108; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=111; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=124; I='ba'; J='Appreciation'},
    @{Row=165; I='sv'; J='Statement-opinion'},
    @{Row=170; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=171; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=174; I='sd'; J='Statement-non-opinion'},
    @{Row=184; I='sd'; J='Statement-non-opinion'},
    @{Row=193; I='sd'; J='Statement-non-opinion'},
    @{Row=204; I='sv'; J='Statement-opinion'},
    @{Row=217; I='%'; J='Uninterpretable'},
    @{Row=227; I='sv'; J='Statement-opinion'},
    @{Row=254; I='ba'; J='Appreciation'},
    @{Row=266; I='ba'; J='Appreciation'},
    @{Row=268; I='ba'; J='Appreciation'},
    @{Row=273; I='ba'; J='Appreciation'},
    @{Row=274; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=287; I='ba'; J='Appreciation'},
    @{Row=294; I='sd'; J='Statement-non-opinion'},
    @{Row=299; I='sd'; J='Statement-non-opinion'},
    @{Row=316; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=328; I='ba'; J='Appreciation'},
    @{Row=350; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=352; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=366; I='aa'; J='Agree/Accept'},
    @{Row=368; I='%'; J='Uninterpretable'},
    @{Row=371; I='sv'; J='Statement-opinion'},
    @{Row=399; I='aa'; J='Agree/Accept'},
    @{Row=419; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=466; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=470; I='sd'; J='Statement-non-opinion'},
    @{Row=471; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=474; I='sd'; J='Statement-non-opinion'},
    @{Row=505; I='aa'; J='Agree/Accept'},
    @{Row=534; I='sd'; J='Statement-non-opinion'},
    @{Row=538; I='ba'; J='Appreciation'},
    @{Row=544; I='ba'; J='Appreciation'},
    @{Row=550; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=551; I='sd'; J='Statement-non-opinion'},
    @{Row=553; I='sd'; J='Statement-non-opinion'},
    @{Row=555; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=591; I='sv'; J='Statement-opinion'},
    @{Row=593; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"
